$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the default/original style from an untouched header-adjacent data cell
$origStyle = $ws.Range("B2").Style

# Cells whose new price text would otherwise be auto-parsed as a number by Excel;
# temporarily mark them as Text format so the literal string is preserved.
$textForceRows = @(5,6,7,9,10,11,19,23,24,25,26,28,29,30,32,34,38,40,43,47,48,51)
foreach ($rn in $textForceRows) {
    $ws.Range("D$rn").NumberFormat = "@"
}

# Apply the updated values from the crypto price/volume refresh
$ws.Range("D2").Value = "49.865.57"
$ws.Range("E2").Value = "  +3.71%  "
$ws.Range("D3").Value = "2.644.68"
$ws.Range("E3").Value = "  +5.91%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "113.77"
$ws.Range("E5").Value = "  +7.58%  "
$ws.Range("D6").Value = "326.49"
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("D7").Value = "0.530"
$ws.Range("E7").Value = "  +1.64%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.553"
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("D10").Value = "41.07"
$ws.Range("E10").Value = "  +5.78%  "
$ws.Range("D11").Value = "20.16"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("E12").Value = "  +2.24%  "
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("E14").Value = "  +3.82%  "
$ws.Range("D15").Value = "3.059.97"
$ws.Range("E15").Value = "  +5.97%  "
$ws.Range("D16").Value = "2.640.83"
$ws.Range("E17").Value = "  +4.68%  "
$ws.Range("D18").Value = "49.780.35"
$ws.Range("E18").Value = "  +3.89%  "
$ws.Range("D19").Value = "13.14"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").Value = "0.0₃0957"
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("D23").Value = "72.04"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("D24").Value = "277.25"
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("D26").Value = "26.75"
$ws.Range("E26").Value = "  +3.81%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "9.98"
$ws.Range("E28").Value = "  +2.73%  "
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  -2.68%  "
$ws.Range("D30").Value = "36.03"
$ws.Range("E30").Value = "  +3.70%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "50.10"
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("E33").Value = "  +2.71%  "
$ws.Range("D34").Value = "19.50"
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("E35").Value = "  +4.38%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +7.02%  "
$ws.Range("D38").Value = "4.85"
$ws.Range("E38").Value = "  +6.00%  "
$ws.Range("E39").Value = "  +8.25%  "
$ws.Range("D40").Value = "126.66"
$ws.Range("E40").Value = "  +3.76%  "
$ws.Range("E41").Value = "  +1.84%  "
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("D43").Value = "22.02"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("E44").Value = "  +3.74%  "
$ws.Range("D45").Value = "2.081.13"
$ws.Range("E45").Value = "  +4.11%  "
$ws.Range("E46").Value = "  +5.79%  "
$ws.Range("D47").Value = "2.31"
$ws.Range("E47").Value = "  +15.44%  "
$ws.Range("D48").Value = "1.97"
$ws.Range("E48").Value = "  +4.62%  "
$ws.Range("E49").Value = "  +2.36%  "
$ws.Range("E50").Value = "  +4.25%  "
$ws.Range("D51").Value = "60.07"
$ws.Range("E51").Value = "  +6.11%  "

# Restore original (default) style on the temporarily reformatted cells
foreach ($rn in $textForceRows) {
    $ws.Range("D$rn").Style = $origStyle
}
